$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3413.5454
$ws.Range("I51").Value = 2450
$ws.Range("K51").Value = 2450
$ws.Range("M51").Value = -1966
$ws.Range("H113").Value = 2492.2307
$ws.Range("I113").Value = 2490.818
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 2490.818
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 763.1819999999998
$ws.Range("N113").Value = -9008
$ws.Range("H137").Value = 6452803
$ws.Range("J137").Value = 28573530
$ws.Range("L137").Value = 85720590
$ws.Range("N137").Value = -85725690

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 35717690
$ws.Range("I61").Value = 41670216
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 41670216
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -41670004
$ws.Range("N61").Value = -2924
$ws.Range("H63").Value = 2710.625
$ws.Range("I63").Value = 2710.625
$ws.Range("K63").Value = 2710.625
$ws.Range("M63").Value = -2024.625
$ws.Range("H66").Value = 2710.625
$ws.Range("I66").Value = 2710.625
$ws.Range("K66").Value = 13553.125
$ws.Range("M66").Value = -10121.125
$ws.Range("H74").Value = 12823003
$ws.Range("I74").Value = 17242520
$ws.Range("K74").Value = 17242520
$ws.Range("M74").Value = -17241646
$ws.Range("H77").Value = 12823003
$ws.Range("I77").Value = 17242520
$ws.Range("K77").Value = 86212600
$ws.Range("M77").Value = -86208232
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 3135.5454
$ws.Range("I132").Value = 2515.8823
$ws.Range("J132").Value = 5242.4
$ws.Range("K132").Value = 7547.646900000001
$ws.Range("L132").Value = 15727.2
$ws.Range("M132").Value = -5017.646900000001
$ws.Range("N132").Value = -20787.2
$ws.Range("H136").Value = 35717690
$ws.Range("I136").Value = 41670216
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 125010648
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -125008098
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 27684.285
$ws.Range("J35").Value = 27684.285
$ws.Range("L35").Value = 27684.285
$ws.Range("N35").Value = -28304.285
$ws.Range("H53").Value = 47820
$ws.Range("J53").Value = 47820
$ws.Range("L53").Value = 47820
$ws.Range("N53").Value = -48968
$ws.Range("H82").Value = 17672.6
$ws.Range("I82").Value = 17757
$ws.Range("J82").Value = 17616.334
$ws.Range("K82").Value = 17757
$ws.Range("L82").Value = 17616.334
$ws.Range("M82").Value = -17374
$ws.Range("N82").Value = -18382.334
$ws.Range("H85").Value = 17672.6
$ws.Range("I85").Value = 17757
$ws.Range("J85").Value = 17616.334
$ws.Range("K85").Value = 17757
$ws.Range("L85").Value = 17616.334
$ws.Range("M85").Value = -16431
$ws.Range("N85").Value = -20268.334
$ws.Range("H94").Value = 1709.1428
$ws.Range("I94").Value = 1432.5714
$ws.Range("J94").Value = 1985.7142
$ws.Range("K94").Value = 1432.5714
$ws.Range("L94").Value = 1985.7142
$ws.Range("M94").Value = -981.5714
$ws.Range("N94").Value = -2887.7142
$ws.Range("H105").Value = 4535.7417
$ws.Range("I105").Value = 3651.5
$ws.Range("K105").Value = 3651.5
$ws.Range("M105").Value = -1904.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9808126
$ws.Range("I31").Value = 4962.963
$ws.Range("J31").Value = 47620330
$ws.Range("K31").Value = 4962.963
$ws.Range("L31").Value = 47620330
$ws.Range("M31").Value = -4667.963
$ws.Range("N31").Value = -47620920
$ws.Range("H34").Value = 9808126
$ws.Range("I34").Value = 4962.963
$ws.Range("J34").Value = 47620330
$ws.Range("K34").Value = 4962.963
$ws.Range("L34").Value = 47620330
$ws.Range("M34").Value = -4760.963
$ws.Range("N34").Value = -47620734
$ws.Range("H50").Value = 7973.6
$ws.Range("J50").Value = 7973.6
$ws.Range("L50").Value = 7973.6
$ws.Range("N50").Value = -9223.6
$ws.Range("H59").Value = 22743.666
$ws.Range("I59").Value = 104
$ws.Range("J59").Value = 34063.5
$ws.Range("K59").Value = 104
$ws.Range("L59").Value = 34063.5
$ws.Range("M59").Value = 1041
$ws.Range("N59").Value = -36353.5
$ws.Range("H60").Value = 7930.2856
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 7585.3335
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 7585.3335
$ws.Range("M60").Value = -9489
$ws.Range("N60").Value = -8607.333500000001
$ws.Range("H68").Value = 20295
$ws.Range("J68").Value = 20295
$ws.Range("L68").Value = 20295
$ws.Range("N68").Value = -21793
$ws.Range("H71").Value = 20295
$ws.Range("J71").Value = 20295
$ws.Range("L71").Value = 60885
$ws.Range("N71").Value = -68373
$ws.Range("H74").Value = 38314
$ws.Range("J74").Value = 38314
$ws.Range("L74").Value = 38314
$ws.Range("N74").Value = -40062
$ws.Range("H77").Value = 38314
$ws.Range("J77").Value = 38314
$ws.Range("L77").Value = 114942
$ws.Range("N77").Value = -123678
$ws.Range("H105").Value = 2750
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 2750
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 2750
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -6244
$ws.Range("H132").Value = 41678520
$ws.Range("I132").Value = 55568916
$ws.Range("J132").Value = 7337.3335
$ws.Range("K132").Value = 166706748
$ws.Range("L132").Value = 22012.0005
$ws.Range("M132").Value = -166704218
$ws.Range("N132").Value = -27072.0005
$ws.Range("H140").Value = 28009.092
$ws.Range("J140").Value = 28009.092
$ws.Range("L140").Value = 28009.092
$ws.Range("N140").Value = -38369.092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 3288.8
$ws.Range("J88").Value = 3288.8
$ws.Range("L88").Value = 9866.400000000001
$ws.Range("N88").Value = -10722.4
$ws.Range("H91").Value = 3288.8
$ws.Range("J91").Value = 3288.8
$ws.Range("L91").Value = 9866.400000000001
$ws.Range("N91").Value = -12830.4
$ws.Range("H121").Value = 1006.2105
$ws.Range("I121").Value = 500
$ws.Range("J121").Value = 1101.125
$ws.Range("K121").Value = 1500
$ws.Range("L121").Value = 3303.375
$ws.Range("M121").Value = -190
$ws.Range("N121").Value = -5923.375
$ws.Range("H122").Value = 1459.6
$ws.Range("I122").Value = 1201.5333
$ws.Range("J122").Value = 1717.6666
$ws.Range("K122").Value = 10813.7997
$ws.Range("L122").Value = 15458.9994
$ws.Range("M122").Value = -8363.7997
$ws.Range("N122").Value = -20358.9994
$ws.Range("H134").Value = 4312.609
$ws.Range("I134").Value = 2842.1428
$ws.Range("J134").Value = 6600
$ws.Range("K134").Value = 8526.428400000001
$ws.Range("L134").Value = 19800
$ws.Range("M134").Value = -3456.428400000001
$ws.Range("N134").Value = -29940
$ws.Range("H137").Value = 6518.7856
$ws.Range("I137").Value = 2682.5
$ws.Range("J137").Value = 8053.3
$ws.Range("K137").Value = 8047.5
$ws.Range("L137").Value = 24159.9
$ws.Range("M137").Value = -2947.5
$ws.Range("N137").Value = -34359.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4144.8
$ws.Range("I126").Value = 2562.25
$ws.Range("J126").Value = 5199.8335
$ws.Range("K126").Value = 7686.75
$ws.Range("L126").Value = 15599.5005
$ws.Range("M126").Value = -5216.75
$ws.Range("N126").Value = -20539.5005
$ws.Range("H132").Value = 5825.92
$ws.Range("I132").Value = 4783.8184
$ws.Range("J132").Value = 6644.7144
$ws.Range("K132").Value = 14351.4552
$ws.Range("L132").Value = 19934.1432
$ws.Range("M132").Value = -11821.4552
$ws.Range("N132").Value = -24994.1432
$ws.Range("H138").Value = 58539.4
$ws.Range("J138").Value = 58539.4
$ws.Range("L138").Value = 58539.4
$ws.Range("N138").Value = -68819.39999999999
$ws.Range("H140").Value = 50950
$ws.Range("J140").Value = 50950
$ws.Range("L140").Value = 50950
$ws.Range("N140").Value = -61310

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H81").Value = 39268.6
$ws.Range("J81").Value = 39268.6
$ws.Range("L81").Value = 39268.6
$ws.Range("N81").Value = -41264.6
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H84").Value = 39268.6
$ws.Range("J84").Value = 39268.6
$ws.Range("L84").Value = 117805.8
$ws.Range("N84").Value = -127789.8
$ws.Range("H122").Value = 6282.7
$ws.Range("I122").Value = 9035.571
$ws.Range("K122").Value = 27106.713
$ws.Range("M122").Value = -24656.713
$ws.Range("H132").Value = 14295020
$ws.Range("I132").Value = 7028.143
$ws.Range("J132").Value = 23820348
$ws.Range("K132").Value = 21084.429
$ws.Range("L132").Value = 71461044
$ws.Range("M132").Value = -18554.429
$ws.Range("N132").Value = -71466104

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2947.8948
$ws.Range("I132").Value = 2026.9
$ws.Range("J132").Value = 3971.2222
$ws.Range("K132").Value = 6080.700000000001
$ws.Range("L132").Value = 11913.6666
$ws.Range("M132").Value = -3550.700000000001
$ws.Range("N132").Value = -16973.6666
$ws.Range("H136").Value = 1446.6086
$ws.Range("I136").Value = 1322.1765
$ws.Range("J136").Value = 1799.1666
$ws.Range("K136").Value = 3966.5295
$ws.Range("L136").Value = 5397.4998
$ws.Range("M136").Value = -1416.5295
$ws.Range("N136").Value = -10497.4998
$ws.Range("H138").Value = 65424.5
$ws.Range("J138").Value = 65424.5
$ws.Range("L138").Value = 65424.5
